# ---------------------------------------------------------------------------
# "added L3 to course Mapp"
#
# 1. Insert a new "Foundation" (L3) worksheet at the front of the workbook,
#    listing the Computing Foundation modules.
# 2. Rename "UG Map" -> "UG Map L4-L6" (now that L3 has its own sheet).
# 3. Leave "PG Map" / "PG MAIDS" as-is, but make "PG MAIDS" the active tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: convert a "stored" OOXML column-width (as it would be written to
# <col width="...">) into the value to feed Excel's ColumnWidth COM property
# (character units get padded by the grid-line pixels, so the two scales
# disagree). Using MaximumDigitWidth = 6px, as Excel computes internally.
function StoredWidthToColumnWidth([double]$target) {
    $mdw = 6
    $pixels = [math]::Floor(((256 * $target + [math]::Floor(128 / $mdw)) / 256) * $mdw)
    return ($pixels - 5) / $mdw
}

# ---------------------------------------------------------------------------
# 1. Add the Foundation sheet and move it to the front.
# ---------------------------------------------------------------------------
$foundation = $wb.Worksheets.Add()
$foundation.Name = "Foundation"
$foundation.Move($wb.Worksheets.Item(1))

$foundation.Columns.Item(1).ColumnWidth = StoredWidthToColumnWidth 6.5
$foundation.Columns.Item(2).ColumnWidth = StoredWidthToColumnWidth 41

# Title block: merged A1:B2 "Computing Foundation"
$title = $foundation.Range("A1:B2")
$title.Merge() | Out-Null
$title.Value = "Computing Foundation"
$title.Font.Bold = $true
$title.Font.Size = 22
$title.HorizontalAlignment = -4108   # xlCenter
$title.VerticalAlignment = -4108     # xlCenter
$title.Interior.Color = 65535        # yellow (matches the other title rows)

$foundation.Rows.Item(1).RowHeight = 29
$foundation.Rows.Item(2).RowHeight = 29

# Thick outer box border around the merged title.
$edges = 7, 8, 9, 10   # xlEdgeLeft, xlEdgeTop, xlEdgeBottom, xlEdgeRight
foreach ($e in $edges) {
    $title.Borders.Item($e).LineStyle = 1
    $title.Borders.Item($e).Weight = 4
}
$title.Borders.Item(7).Weight = 4
$title.Borders.Item(8).Weight = 4
$title.Borders.Item(9).Weight = 2
$title.Borders.Item(10).Weight = 2
$title.Borders.Item(11).LineStyle = 1   # inside vertical
$title.Borders.Item(11).Weight = 2
$title.Borders.Item(12).LineStyle = 1   # inside horizontal
$title.Borders.Item(12).Weight = 2

# Module rows (1-6), same look as the numbered lists on the other sheets:
# column A = index number, column B = module name; banded row shading,
# thin borders, medium border down the outer (left/right) edges.
$modules = @(
    "COM300 Problem Solving",
    "COM302 Group Technology Project",
    "COM304 Foundation Computing",
    "COM305 Induviudual Degree Related Project",
    "COM306 Digital Media Technologies",
    "COM307 Foundation Mathematics"
)

for ($i = 0; $i -lt $modules.Length; $i++) {
    $r = 3 + $i
    $num = $i + 1
    $shaded = ($num % 2) -eq 1

    $numCell = $foundation.Range("A$r")
    $numCell.Value = $num
    $numCell.Font.Bold = $true
    $numCell.HorizontalAlignment = -4108
    $numCell.VerticalAlignment = -4108
    if ($shaded) {
        $numCell.Interior.Color = 15000804
    }

    $txtCell = $foundation.Range("B$r")
    $txtCell.Value = $modules[$i]
    if ($shaded) {
        $txtCell.Interior.Color = 15000804
    }

    foreach ($c in @($numCell, $txtCell)) {
        $c.Borders.Item(7).LineStyle = 1   # left
        $c.Borders.Item(8).LineStyle = 1   # top
        $c.Borders.Item(9).LineStyle = 1   # bottom
        $c.Borders.Item(10).LineStyle = 1  # right
        $c.Borders.Item(7).Weight = 2
        $c.Borders.Item(8).Weight = 2
        $c.Borders.Item(9).Weight = 2
        $c.Borders.Item(10).Weight = 2
    }
    $numCell.Borders.Item(7).Weight = 4   # outer-left thicker
    $txtCell.Borders.Item(10).Weight = 4  # outer-right thicker
}

# Bottom row of the table gets a thicker closing border underneath.
$lastA = $foundation.Range("A8")
$lastB = $foundation.Range("B8")
$lastA.Borders.Item(9).Weight = 4
$lastB.Borders.Item(9).Weight = 4

$foundation.Range("F8").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Rename "UG Map" to "UG Map L4-L6".
# ---------------------------------------------------------------------------
$ugMap = $wb.Worksheets.Item("UG Map")
$ugMap.Name = "UG Map L4-L6"
$ugMap.Range("B41").Select() | Out-Null

$pgMap = $wb.Worksheets.Item("PG Map")
$pgMap.Range("D28").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Make "PG MAIDS" the active sheet/tab.
# ---------------------------------------------------------------------------
$pgMaids = $wb.Worksheets.Item("PG MAIDS")
$pgMaids.Select()
$pgMaids.Range("B41").Select() | Out-Null
